$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the date-cell number format used throughout column D before the insert
$dateFmt = $ws.Range("D349").NumberFormat

# Insert a new row at position 348, shifting existing rows 348-398 down to 349-399
$ws.Rows.Item(348).Insert()

# Re-apply the date number format to the new D348 cell (Insert does not always carry it)
$ws.Range("D348").NumberFormat = $dateFmt

# Populate the newly inserted row 348 with the new record
$ws.Range("A348").Value = 10
$ws.Range("B348").Value = "Vega Modelo de Temuco"
$ws.Range("C348").Value = "La Araucanía"
$ws.Range("D348").Value = 44491
$ws.Range("E348").Value = 9
$ws.Range("F348").Value = 100112003
$ws.Range("G348").Value = "Ajo"
$ws.Range("H348").Value = "Chino"
$ws.Range("I348").Value = "Primera"
$ws.Range("J348").Value = 220
$ws.Range("K348").Value = 19000
$ws.Range("L348").Value = 20000
$ws.Range("M348").Value = 19295
$ws.Range("N348").Value = "$/malla 10 kilos"
$ws.Range("O348").Value = "China"
$ws.Range("P348").Value = 1930
$ws.Range("Q348").Value = 10
$ws.Range("R348").Value = "Hortaliza"
